$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and two pairs of swapped rows)

# Row 2: Price, Volume
$ws.Range("D2").Value = "43.033.48"
$ws.Range("E2").Value = "  +0.61%  "

# Row 3: Price
$ws.Range("D3").Value = "2.366.44"

# Row 4: Price, Volume
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5: Price, Volume
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.45%  "

# Row 6: Price, Volume
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "

# Row 7: Volume
$ws.Range("E7").Value = "  -0.22%  "

# Row 8: Volume
$ws.Range("E8").Value = "  -0.02%  "

# Row 9: Volume
$ws.Range("E9").Value = "  -0.23%  "

# Row 10: Price, Volume
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.11%  "

# Row 11: Volume
$ws.Range("E11").Value = "  +3.62%  "

# Row 13: Price, Volume
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.38"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.05%  "

# Row 14: Volume
$ws.Range("E14").Value = "  +0.44%  "

# Row 15: Price, Volume
$ws.Range("D15").Value = "2.732.33"
$ws.Range("E15").Value = "  +2.25%  "

# Row 16: Price, Volume
$ws.Range("D16").Value = "2.369.48"
$ws.Range("E16").Value = "  +2.18%  "

# Row 17: Volume
$ws.Range("E17").Value = "  +0.90%  "

# Row 18: Price, Volume
$ws.Range("D18").Value = "42.999.75"
$ws.Range("E18").Value = "  +0.65%  "

# Row 19: Price, Volume
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.23%  "

# Row 20: Price, Volume
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.25%  "

# Row 21: Volume
$ws.Range("E21").Value = "  -0.39%  "

# Row 22: Price, Volume
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "

# Row 23: Price, Volume
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.08%  "

# Row 24: Price, Volume
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.51%  "

# Row 25: Coin, Link, Price, Volume
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.86%  "

# Row 26: Coin, Link, Price, Volume
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "

# Row 27: Price, Volume
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.15%  "

# Row 28: Volume
$ws.Range("E28").Value = "  +0.34%  "

# Row 29: Price, Volume
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.68%  "

# Row 30: Price, Volume
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.73%  "

# Row 31: Volume
$ws.Range("E31").Value = "  -0.06%  "

# Row 32: Price, Volume
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.43%  "

# Row 33: Price, Volume
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.14%  "

# Row 34: Price, Volume
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0720"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.22%  "

# Row 35: Volume
$ws.Range("E35").Value = "  +4.25%  "

# Row 36: Volume
$ws.Range("E36").Value = "  +3.68%  "

# Row 37: Volume
$ws.Range("E37").Value = "  -2.30%  "

# Row 38: Volume
$ws.Range("E38").Value = "  -1.64%  "

# Row 39: Coin, Link, Price, Volume
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.80"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.23%  "

# Row 40: Coin, Link, Price, Volume
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "122.65"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -26.23%  "

# Row 41: Volume
$ws.Range("E41").Value = "  -0.50%  "

# Row 42: Price, Volume
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.26%  "

# Row 43: Price, Volume
$ws.Range("D43").Value = "1.937.51"
$ws.Range("E43").Value = "  +0.63%  "

# Row 44: Volume
$ws.Range("E44").Value = "  +0.25%  "

# Row 45: Price, Volume
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.14"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.42%  "

# Row 46: Volume
$ws.Range("E46").Value = "  -0.07%  "

# Row 47: Price, Volume
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -8.27%  "

# Row 48: Price
$ws.Range("D48").Value = "2.590.05"

# Row 49: Volume
$ws.Range("E49").Value = "  +2.46%  "

# Row 50: Coin, Link, Price, Volume
$ws.Range("B50").Value = "TrustWalletToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.78%  "

# Row 51: Coin, Link, Price, Volume
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.00%  "

